# Mini cart validation for osprey
# Adds a new "Minicart" worksheet (after the last existing sheet) that
# mirrors the other per-flow validation sheets (DataSet, Header, ...):
# a yellow header row followed by the captured account / product / quantity
# data for the mini-cart flow.

$wb = $excel.ActiveWorkbook

# --- Add the new sheet at the very end of the tab strip -------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Minicart"

# --- Header row (row 1) : yellow-filled labels -----------------------------
$headers = @(
    "DataSet", "UserName", "Prod UserName", "Email", "Password",
    "Confirm Password", "FirstName", "LastName", "Street", "City",
    "Region", "postcode", "phone", "Products", "Color", "Size",
    "Colorproduct", "Quantity", "methods", "cardNumber", "ExpMonthYear", "cvv"
)
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $ws.Cells.Item(1, $i + 1)
    $cell.Value = $headers[$i]
    $cell.Interior.Color = 65535
}

# --- Row 2 : account / address details -------------------------------------
$ws.Range("A2").Value = "Account"
$ws.Range("B2").Value = "testersemail.278@gmail.com"
$ws.Range("D2").Value = "testersemail.278@gmail.com"
$ws.Range("E2").Value = "Testers@278"
$ws.Range("F2").Value = "Testers@278"
$ws.Range("G2").Value = "QA"
$ws.Range("H2").Value = "TEST"
$ws.Range("I2").Value = "6 Walnut Valley Dr"
$ws.Range("J2").Value = "Little Rock"
$ws.Range("K2").Value = "Arkansas"
$ws.Range("L2").Value = 72211
$ws.Range("M2").Value = 9898989898

# --- Row 3 : simple product ------------------------------------------------
$ws.Range("A3").Value = "Product"
$ws.Range("N3").Value = "POCO® CARRYING CASE"
$ws.Range("O3").Value = "Black"
$ws.Range("P3").Value = "S/M"
$ws.Range("Q3").Value = "AETHER™ 55"
$ws.Range("R3").Value = "'1"

# --- Row 4 : configurable product ------------------------------------------
$ws.Range("A4").Value = "ConfigurableProduct"
$ws.Range("N4").Value = "AETHER™ 55"
$ws.Range("O4").Value = "Black"
$ws.Range("P4").Value = "S/M"
$ws.Range("R4").Value = "'1"

# --- Row 5 : product quantity -----------------------------------------------
$ws.Range("A5").Value = "Product Qunatity"
$ws.Range("R5").Value = "'2"

# --- Finish with the cursor where the author left it -----------------------
$ws.Range("L11").Select()
